# Update Mappings 22 Ontologies
# Adds a new "SBO_DEF" column (F) to Sheet1, populated with the literal
# string "[]" for every existing data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column F, matching the formatting of the other header
# cells (B1:E1): bold font, centered/top aligned, thin box border.
$ws.Range("F1").Value = "SBO_DEF"
$ws.Range("F1").Font.Bold = $true
$ws.Range("F1").HorizontalAlignment = -4108
$ws.Range("F1").VerticalAlignment = -4160
$ws.Range("F1").Borders.LineStyle = 1

# Find the last used data row (header is row 1, data starts row 2)
$lastRow = $ws.Cells.Item($ws.Rows.Count, "B").End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 6).Value = "[]"
}
